$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    # Write a value while forcing the cell to remain text-typed (no
    # numeric auto-coercion of numeric-looking strings) and without
    # leaving a residual style/number-format behind.
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.ClearFormats()
}

Set-TextValue "D2" "61.924.97"
Set-TextValue "E2" "  -1.60%  "

Set-TextValue "D3" "2.915.97"
Set-TextValue "E3" "  -2.13%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.06%  "

Set-TextValue "D5" "587.49"
Set-TextValue "E5" "  -1.60%  "

Set-TextValue "D6" "147.03"
Set-TextValue "E6" "  +1.26%  "

Set-TextValue "E7" "  +0.02%  "

Set-TextValue "D8" "0.507"
Set-TextValue "E8" "  +0.59%  "

Set-TextValue "D9" "2.916.51"
Set-TextValue "E9" "  -2.10%  "

Set-TextValue "D10" "7.06"
Set-TextValue "E10" "  -3.74%  "

Set-TextValue "E11" "  +4.36%  "

Set-TextValue "E12" "  -2.92%  "

Set-TextValue "D13" "0.0000240"
Set-TextValue "E13" "  +1.53%  "

Set-TextValue "D14" "33.17"
Set-TextValue "E14" "  -1.10%  "

Set-TextValue "E15" "  -1.57%  "

Set-TextValue "D16" "3.397.02"
Set-TextValue "E16" "  -2.15%  "

Set-TextValue "D17" "61.939.94"
Set-TextValue "E17" "  -1.36%  "

Set-TextValue "E18" "  -1.90%  "

Set-TextValue "D19" "2.917.24"
Set-TextValue "E19" "  -1.92%  "

Set-TextValue "D20" "435.72"
Set-TextValue "E20" "  -1.52%  "

Set-TextValue "D21" "13.46"
Set-TextValue "E21" "  -0.92%  "

Set-TextValue "E22" "  -2.81%  "

Set-TextValue "D23" "6.96"
Set-TextValue "E23" "  -2.83%  "

Set-TextValue "D24" "81.15"
Set-TextValue "E24" "  -1.39%  "

Set-TextValue "E25" "  -1.90%  "

Set-TextValue "D26" "10.26"
Set-TextValue "E26" "  -5.67%  "

Set-TextValue "E28" "  +0.02%  "

Set-TextValue "E29" "  +22.64%  "

Set-TextValue "D30" "7.29"
Set-TextValue "E30" "  +2.69%  "

Set-TextValue "D31" "2.57"
Set-TextValue "E31" "  -2.08%  "

Set-TextValue "D32" "2.12"
Set-TextValue "E32" "  -0.61%  "

Set-TextValue "E33" "  +1.36%  "

Set-TextValue "D34" "25.97"
Set-TextValue "E34" "  -2.29%  "

Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  -0.13%  "

Set-TextValue "E36" "  -1.41%  "

Set-TextValue "D37" "3.09"
Set-TextValue "E37" "  +4.39%  "

Set-TextValue "E38" "  -2.42%  "

Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "2.02"
Set-TextValue "E39" "  -0.55%  "

Set-TextValue "B40" "OKB"
Set-TextValue "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D40" "49.18"
Set-TextValue "E40" "  -0.89%  "

Set-TextValue "D41" "8.37"
Set-TextValue "E41" "  -2.97%  "

Set-TextValue "D42" "0.116"
Set-TextValue "E42" "  -2.02%  "

Set-TextValue "D43" "0.274"
Set-TextValue "E43" "  -4.21%  "

Set-TextValue "D44" "39.18"
Set-TextValue "E44" "  +1.15%  "

Set-TextValue "D45" "2.701.38"
Set-TextValue "E45" "  -0.20%  "

Set-TextValue "D46" "134.48"
Set-TextValue "E46" "  +0.24%  "

Set-TextValue "D48" "345.92"
Set-TextValue "E48" "  -6.71%  "

Set-TextValue "E49" "  +0.05%  "

Set-TextValue "E50" "  -1.31%  "

Set-TextValue "D51" "22.42"
Set-TextValue "E51" "  -4.08%  "
